$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update masthead text (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# --- Simple numeric value updates (style/number-format unchanged) ---
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 50
$ws.Range("L14").Value = -25
$ws.Range("N14").Value = 0
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -60
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 500
$ws.Range("F16").Value = 16
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 42
$ws.Range("J16").Value = 39
$ws.Range("K16").Value = 7.692307692307
$ws.Range("L16").Value = 180
$ws.Range("M16").Value = -19.230769230769
$ws.Range("N16").Value = -77.419354838709
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -19.047619047619
$ws.Range("I17").Value = 70
$ws.Range("J17").Value = 84
$ws.Range("K17").Value = -16.666666666666
$ws.Range("L17").Value = 48.936170212766
$ws.Range("M17").Value = 59.090909090909
$ws.Range("N17").Value = 0
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 17
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -41.379310344827
$ws.Range("L18").Value = -19.047619047619
$ws.Range("M18").Value = -77.027027027027
$ws.Range("N18").Value = -91.707317073170
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -71.428571428571
$ws.Range("F19").Value = 13
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = -18.75
$ws.Range("I19").Value = 66
$ws.Range("J19").Value = 79
$ws.Range("K19").Value = -16.455696202531
$ws.Range("L19").Value = 112.903225806452
$ws.Range("M19").Value = 1.538461538461
$ws.Range("N19").Value = -39.449541284403
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 500
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 8
$ws.Range("L20").Value = 125
$ws.Range("M20").Value = 58.823529411764
$ws.Range("N20").Value = -92.252510760401
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 68
$ws.Range("H21").Value = 1.470588235294
$ws.Range("I21").Value = 254
$ws.Range("J21").Value = 288
$ws.Range("K21").Value = -11.805555555555
$ws.Range("L21").Value = 72.789115646258
$ws.Range("M21").Value = -6.273062730627
$ws.Range("N21").Value = -80.109631949882
$ws.Range("L22").Value = -40
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = 14.285714285714
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = 30.769230769230
$ws.Range("L23").Value = 36
$ws.Range("M23").Value = 183.333333333333
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 46
$ws.Range("G24").Value = 45
$ws.Range("H24").Value = 2.222222222222
$ws.Range("I24").Value = 198
$ws.Range("J24").Value = 166
$ws.Range("K24").Value = 19.277108433734
$ws.Range("L24").Value = 62.295081967213
$ws.Range("M24").Value = 44.525547445255
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = -32.432432432432
$ws.Range("I25").Value = 108
$ws.Range("J25").Value = 111
$ws.Range("K25").Value = -2.702702702702
$ws.Range("L25").Value = 33.333333333333
$ws.Range("M25").Value = -11.475409836065
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -100
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = -63.636363636363
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = -12.5
$ws.Range("L27").Value = 0
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 13
$ws.Range("K28").Value = -46.153846153846
$ws.Range("L28").Value = -41.666666666666
$ws.Range("M28").Value = 40
$ws.Range("N28").Value = -30
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 7
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = -30
$ws.Range("L29").Value = -22.222222222222
$ws.Range("M29").Value = 75
$ws.Range("N29").Value = -22.222222222222

# --- Cells that change from text ("0"/"***.*") to numeric (apply number format + value) ---
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 2
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = -100
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 2
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 2
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C29").Value = 2

# --- Cells that change from numeric back to text (copy format from a General-formatted text cell, then set text value) ---
$generalTextSource = $ws.Range("C15")
$ws.Range("F26").Value = "'0"
$generalTextSource.Copy()
$ws.Range("F26").PasteSpecial(-4122)

$excel.CutCopyMode = $false
